$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$src = $ws.Range("D2")
$dst = $ws.Range("D60")
$dst.Value = 0
$src.Copy()
$dst.PasteSpecial(-4122)
$dst.Borders.Item(8).LineStyle = -4142
$dst.Borders.Item(9).LineStyle = -4142
